# Auto-generated Excel COM-interop script applying the KCOR_summary.xlsx diff.
# For each sheet: update the KCOR / CI_Lower / CI_Upper text figures in place.
# Numeric-looking values are written with a leading apostrophe so Excel keeps
# storing them as text (matching the workbook's inlineStr convention) instead
# of silently converting them to numbers.

$wb = $excel.ActiveWorkbook

# ---- sheet1 (Worksheets.Item(1)) ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 3).Value = "'1.5012"
$ws.Cells.Item(4, 4).Value = "'1.459"
$ws.Cells.Item(4, 5).Value = "'1.545"
$ws.Cells.Item(5, 3).Value = "'1.1943"
$ws.Cells.Item(5, 4).Value = "'1.102"
$ws.Cells.Item(5, 5).Value = "'1.294"
$ws.Cells.Item(6, 3).Value = "'1.2437"
$ws.Cells.Item(6, 4).Value = "'1.196"
$ws.Cells.Item(6, 5).Value = "'1.293"
$ws.Cells.Item(7, 3).Value = "'1.8376"
$ws.Cells.Item(7, 4).Value = "'1.773"
$ws.Cells.Item(7, 5).Value = "'1.905"
$ws.Cells.Item(8, 3).Value = "'1.5456"
$ws.Cells.Item(8, 4).Value = "'1.419"
$ws.Cells.Item(8, 5).Value = "'1.683"
$ws.Cells.Item(10, 3).Value = "'1.5770"
$ws.Cells.Item(10, 4).Value = "'1.173"
$ws.Cells.Item(10, 5).Value = "'2.120"
$ws.Cells.Item(11, 3).Value = "'1.6376"
$ws.Cells.Item(11, 4).Value = "'0.876"
$ws.Cells.Item(11, 5).Value = "'3.062"
$ws.Cells.Item(13, 3).Value = "'8.0955"
$ws.Cells.Item(13, 5).Value = "'80.955"
$ws.Cells.Item(16, 3).Value = "'1.6843"
$ws.Cells.Item(16, 4).Value = "'1.629"
$ws.Cells.Item(16, 5).Value = "'1.741"
$ws.Cells.Item(17, 3).Value = "'1.5732"
$ws.Cells.Item(17, 4).Value = "'1.478"
$ws.Cells.Item(17, 5).Value = "'1.675"
$ws.Cells.Item(18, 3).Value = "'1.9558"
$ws.Cells.Item(18, 4).Value = "'1.886"
$ws.Cells.Item(18, 5).Value = "'2.029"
$ws.Cells.Item(19, 3).Value = "'1.8148"
$ws.Cells.Item(19, 4).Value = "'1.717"
$ws.Cells.Item(19, 5).Value = "'1.919"
$ws.Cells.Item(20, 3).Value = "'1.3536"
$ws.Cells.Item(20, 4).Value = "'1.220"
$ws.Cells.Item(20, 5).Value = "'1.502"
$ws.Cells.Item(22, 3).Value = "'1.3922"
$ws.Cells.Item(22, 4).Value = "'0.996"
$ws.Cells.Item(22, 5).Value = "'1.945"
$ws.Cells.Item(23, 3).Value = "'0.6550"
$ws.Cells.Item(23, 4).Value = "'0.438"
$ws.Cells.Item(23, 5).Value = "'0.979"
$ws.Cells.Item(25, 3).Value = "'3.2189"
$ws.Cells.Item(25, 5).Value = "'32.189"
$ws.Cells.Item(28, 3).Value = "'1.1203"
$ws.Cells.Item(28, 4).Value = "'1.075"
$ws.Cells.Item(28, 5).Value = "'1.168"

# ---- sheet2 (Worksheets.Item(2)) ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 3).Value = "'1.1956"
$ws.Cells.Item(4, 4).Value = "'1.161"
$ws.Cells.Item(4, 5).Value = "'1.231"
$ws.Cells.Item(5, 3).Value = "'0.9462"
$ws.Cells.Item(5, 4).Value = "'0.855"
$ws.Cells.Item(5, 5).Value = "'1.047"
$ws.Cells.Item(6, 3).Value = "'1.0840"
$ws.Cells.Item(6, 4).Value = "'1.025"
$ws.Cells.Item(6, 5).Value = "'1.146"
$ws.Cells.Item(7, 3).Value = "'1.2713"
$ws.Cells.Item(7, 4).Value = "'1.213"
$ws.Cells.Item(7, 5).Value = "'1.332"
$ws.Cells.Item(8, 3).Value = "'1.3714"
$ws.Cells.Item(8, 4).Value = "'1.299"
$ws.Cells.Item(8, 5).Value = "'1.448"
$ws.Cells.Item(9, 3).Value = "'1.3244"
$ws.Cells.Item(9, 4).Value = "'1.220"
$ws.Cells.Item(9, 5).Value = "'1.438"
$ws.Cells.Item(10, 3).Value = "'1.2309"
$ws.Cells.Item(10, 4).Value = "'1.092"
$ws.Cells.Item(10, 5).Value = "'1.388"
$ws.Cells.Item(11, 3).Value = "'0.8595"
$ws.Cells.Item(11, 4).Value = "'0.687"
$ws.Cells.Item(11, 5).Value = "'1.075"
$ws.Cells.Item(13, 3).Value = "'1.0927"
$ws.Cells.Item(13, 4).Value = "'0.443"
$ws.Cells.Item(13, 5).Value = "'2.696"
$ws.Cells.Item(16, 3).Value = "'1.3673"
$ws.Cells.Item(16, 4).Value = "'1.336"
$ws.Cells.Item(16, 5).Value = "'1.399"
$ws.Cells.Item(17, 3).Value = "'1.0828"
$ws.Cells.Item(17, 4).Value = "'1.019"
$ws.Cells.Item(17, 5).Value = "'1.151"
$ws.Cells.Item(18, 3).Value = "'1.2331"
$ws.Cells.Item(18, 4).Value = "'1.189"
$ws.Cells.Item(18, 5).Value = "'1.278"
$ws.Cells.Item(19, 3).Value = "'1.4916"
$ws.Cells.Item(19, 4).Value = "'1.438"
$ws.Cells.Item(19, 5).Value = "'1.547"
$ws.Cells.Item(20, 3).Value = "'1.6234"
$ws.Cells.Item(20, 4).Value = "'1.537"
$ws.Cells.Item(20, 5).Value = "'1.715"
$ws.Cells.Item(21, 3).Value = "'1.3924"
$ws.Cells.Item(21, 4).Value = "'1.255"
$ws.Cells.Item(21, 5).Value = "'1.545"
$ws.Cells.Item(22, 3).Value = "'1.4244"
$ws.Cells.Item(22, 4).Value = "'1.204"
$ws.Cells.Item(22, 5).Value = "'1.685"
$ws.Cells.Item(23, 3).Value = "'1.2366"
$ws.Cells.Item(23, 4).Value = "'0.909"
$ws.Cells.Item(23, 5).Value = "'1.683"
$ws.Cells.Item(25, 3).Value = "'1.1812"
$ws.Cells.Item(25, 4).Value = "'0.479"
$ws.Cells.Item(25, 5).Value = "'2.914"
$ws.Cells.Item(28, 3).Value = "'1.1436"
$ws.Cells.Item(28, 4).Value = "'1.108"
$ws.Cells.Item(28, 5).Value = "'1.180"

# ---- sheet3 (Worksheets.Item(3)) ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 3).Value = "'1.0673"
$ws.Cells.Item(4, 4).Value = "'1.021"
$ws.Cells.Item(4, 5).Value = "'1.115"
$ws.Cells.Item(7, 3).Value = "'0.9575"
$ws.Cells.Item(7, 4).Value = "'0.890"
$ws.Cells.Item(7, 5).Value = "'1.030"
$ws.Cells.Item(10, 3).Value = "'0.9634"
$ws.Cells.Item(10, 4).Value = "'0.812"
$ws.Cells.Item(10, 5).Value = "'1.143"
$ws.Cells.Item(11, 3).Value = "'1.0706"
$ws.Cells.Item(11, 4).Value = "'0.803"
$ws.Cells.Item(11, 5).Value = "'1.427"
$ws.Cells.Item(12, 3).Value = "'1.1050"
$ws.Cells.Item(12, 4).Value = "'0.719"
$ws.Cells.Item(12, 5).Value = "'1.699"
$ws.Cells.Item(16, 3).Value = "'1.0599"
$ws.Cells.Item(16, 4).Value = "'1.036"
$ws.Cells.Item(16, 5).Value = "'1.084"
$ws.Cells.Item(19, 3).Value = "'1.0149"
$ws.Cells.Item(19, 4).Value = "'0.976"
$ws.Cells.Item(19, 5).Value = "'1.056"
$ws.Cells.Item(22, 3).Value = "'1.1585"
$ws.Cells.Item(22, 4).Value = "'1.036"
$ws.Cells.Item(22, 5).Value = "'1.295"
$ws.Cells.Item(23, 3).Value = "'1.1282"
$ws.Cells.Item(23, 4).Value = "'0.933"
$ws.Cells.Item(23, 5).Value = "'1.365"
$ws.Cells.Item(24, 3).Value = "'1.3582"
$ws.Cells.Item(24, 4).Value = "'1.007"
$ws.Cells.Item(24, 5).Value = "'1.832"
$ws.Cells.Item(28, 3).Value = "'0.9931"
$ws.Cells.Item(28, 4).Value = "'0.950"
$ws.Cells.Item(28, 5).Value = "'1.038"
$ws.Cells.Item(40, 3).Value = "'1.6883"
$ws.Cells.Item(40, 4).Value = "'1.655"
$ws.Cells.Item(40, 5).Value = "'1.722"
$ws.Cells.Item(43, 3).Value = "'1.8140"
$ws.Cells.Item(43, 4).Value = "'1.751"
$ws.Cells.Item(43, 5).Value = "'1.879"
$ws.Cells.Item(46, 3).Value = "'1.6409"
$ws.Cells.Item(46, 4).Value = "'1.449"
$ws.Cells.Item(46, 5).Value = "'1.858"
$ws.Cells.Item(47, 3).Value = "'1.1953"
$ws.Cells.Item(47, 4).Value = "'0.971"
$ws.Cells.Item(47, 5).Value = "'1.471"
$ws.Cells.Item(48, 3).Value = "'1.4701"
$ws.Cells.Item(48, 4).Value = "'1.011"
$ws.Cells.Item(48, 5).Value = "'2.139"
$ws.Cells.Item(51, 1).Value = "3 vs 1"
$ws.Cells.Item(52, 3).Value = "'1.5819"
$ws.Cells.Item(52, 4).Value = "'1.515"
$ws.Cells.Item(52, 5).Value = "'1.652"
$ws.Cells.Item(53, 3).Value = "'1.0896"
$ws.Cells.Item(53, 4).Value = "'0.895"
$ws.Cells.Item(53, 5).Value = "'1.326"
$ws.Cells.Item(54, 3).Value = "'1.4574"
$ws.Cells.Item(54, 4).Value = "'1.340"
$ws.Cells.Item(54, 5).Value = "'1.585"
$ws.Cells.Item(55, 3).Value = "'1.8944"
$ws.Cells.Item(55, 4).Value = "'1.764"
$ws.Cells.Item(55, 5).Value = "'2.034"
$ws.Cells.Item(56, 3).Value = "'1.4853"
$ws.Cells.Item(56, 4).Value = "'1.351"
$ws.Cells.Item(56, 5).Value = "'1.633"
$ws.Cells.Item(57, 3).Value = "'1.4670"
$ws.Cells.Item(57, 4).Value = "'1.285"
$ws.Cells.Item(57, 5).Value = "'1.675"
$ws.Cells.Item(58, 3).Value = "'1.7033"
$ws.Cells.Item(58, 4).Value = "'1.413"
$ws.Cells.Item(58, 5).Value = "'2.053"
$ws.Cells.Item(59, 3).Value = "'1.1164"
$ws.Cells.Item(59, 4).Value = "'0.812"
$ws.Cells.Item(59, 5).Value = "'1.536"
$ws.Cells.Item(60, 3).Value = "'1.3305"
$ws.Cells.Item(60, 4).Value = "'0.794"
$ws.Cells.Item(60, 5).Value = "'2.228"
$ws.Cells.Item(61, 3).Value = "'0.6250"
$ws.Cells.Item(61, 4).Value = "'0.226"
$ws.Cells.Item(61, 5).Value = "'1.726"

# ---- sheet4 (Worksheets.Item(4)) ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 3).Value = "'1.3158"
$ws.Cells.Item(4, 4).Value = "'1.276"
$ws.Cells.Item(4, 5).Value = "'1.357"
$ws.Cells.Item(7, 3).Value = "'1.2986"
$ws.Cells.Item(7, 4).Value = "'1.235"
$ws.Cells.Item(7, 5).Value = "'1.365"
$ws.Cells.Item(8, 3).Value = "'1.4239"
$ws.Cells.Item(8, 4).Value = "'1.330"
$ws.Cells.Item(8, 5).Value = "'1.524"
$ws.Cells.Item(10, 3).Value = "'2.0696"
$ws.Cells.Item(10, 4).Value = "'1.503"
$ws.Cells.Item(10, 5).Value = "'2.849"
$ws.Cells.Item(11, 3).Value = "'0.7346"
$ws.Cells.Item(11, 4).Value = "'0.516"
$ws.Cells.Item(11, 5).Value = "'1.045"
$ws.Cells.Item(12, 3).Value = "'0.6613"
$ws.Cells.Item(12, 4).Value = "'0.386"
$ws.Cells.Item(12, 5).Value = "'1.133"
$ws.Cells.Item(13, 3).Value = "'4.2108"
$ws.Cells.Item(13, 5).Value = "'42.108"
$ws.Cells.Item(16, 3).Value = "'1.4395"
$ws.Cells.Item(16, 4).Value = "'1.362"
$ws.Cells.Item(16, 5).Value = "'1.521"
$ws.Cells.Item(28, 3).Value = "'1.2990"
$ws.Cells.Item(28, 4).Value = "'1.258"
$ws.Cells.Item(28, 5).Value = "'1.341"
$ws.Cells.Item(40, 3).Value = "'1.2451"
$ws.Cells.Item(40, 4).Value = "'1.212"
$ws.Cells.Item(40, 5).Value = "'1.279"

# ---- sheet3 (2022_06): append the new '3 vs 2' block (rows 63-74) that holds
#      the figures the old '3 vs 2' row block had before being renamed to '3 vs 1'
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(63, 1).Value = "3 vs 2"
$ws.Cells.Item(64, 2).Value = "ASMR (pooled)"
$ws.Cells.Item(64, 3).Value = "'1.5929"
$ws.Cells.Item(64, 4).Value = "'1.560"
$ws.Cells.Item(64, 5).Value = "'1.626"
$ws.Cells.Item(65, 2).Value = "'1920"
$ws.Cells.Item(65, 3).Value = "'1.4695"
$ws.Cells.Item(65, 4).Value = "'1.369"
$ws.Cells.Item(65, 5).Value = "'1.577"
$ws.Cells.Item(66, 2).Value = "'1930"
$ws.Cells.Item(66, 3).Value = "'1.5195"
$ws.Cells.Item(66, 4).Value = "'1.464"
$ws.Cells.Item(66, 5).Value = "'1.577"
$ws.Cells.Item(67, 2).Value = "'1940"
$ws.Cells.Item(67, 3).Value = "'1.7873"
$ws.Cells.Item(67, 4).Value = "'1.726"
$ws.Cells.Item(67, 5).Value = "'1.851"
$ws.Cells.Item(68, 2).Value = "'1950"
$ws.Cells.Item(68, 3).Value = "'1.5324"
$ws.Cells.Item(68, 4).Value = "'1.458"
$ws.Cells.Item(68, 5).Value = "'1.610"
$ws.Cells.Item(69, 2).Value = "'1960"
$ws.Cells.Item(69, 3).Value = "'1.3517"
$ws.Cells.Item(69, 4).Value = "'1.243"
$ws.Cells.Item(69, 5).Value = "'1.470"
$ws.Cells.Item(70, 2).Value = "'1970"
$ws.Cells.Item(70, 3).Value = "'1.4165"
$ws.Cells.Item(70, 4).Value = "'1.239"
$ws.Cells.Item(70, 5).Value = "'1.619"
$ws.Cells.Item(71, 2).Value = "'1980"
$ws.Cells.Item(71, 3).Value = "'1.0595"
$ws.Cells.Item(71, 4).Value = "'0.837"
$ws.Cells.Item(71, 5).Value = "'1.340"
$ws.Cells.Item(72, 2).Value = "'1990"
$ws.Cells.Item(72, 3).Value = "'1.0824"
$ws.Cells.Item(72, 4).Value = "'0.717"
$ws.Cells.Item(72, 5).Value = "'1.635"
$ws.Cells.Item(73, 2).Value = "'2000"
$ws.Cells.Item(73, 3).Value = "'0.7693"
$ws.Cells.Item(73, 4).Value = "'0.424"
$ws.Cells.Item(73, 5).Value = "'1.397"
